# Updated cryptos list - applying price/volume/coin changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.973.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.44%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.223.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -5.49%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'243.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.34%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -5.53%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'68.13"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -8.21%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.547"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -8.34%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.0956"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.37%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'58.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.61%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'35.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.93%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -2.98%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'6.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.16%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.553.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.59%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'14.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -9.03%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.844"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -6.85%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.235.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.15%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'41.900.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.47%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.0₃0951"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.82%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'72.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -7.34%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -8.36%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'234.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -7.34%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +10.12%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.21%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'3.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.96%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.65%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -2.85%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'9.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.85%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'170.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.33%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'20.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -8.77%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -6.28%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -7.28%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.0712"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.21%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'5.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.16%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'4.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -8.80%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +1.53%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'22.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +17.69%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "'LidoDAOToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'2.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.33%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'VeChain"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.0277"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.32%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'THORChain"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'5.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.37%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'MultiversX"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'66.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.85%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'4.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -10.26%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'8.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.41%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.100"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.84%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.190"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.05%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.13%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'4.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.65%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.23%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'10.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.71%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.79%  "
$ws.Range("E51").Style = "Normal"
